# Applies the "Enigmatic Cosmos / dark energy" -> "Profound Evolution of
# Artistic Expression" rewrite described by the unified diff.
#
# Strategy: the document is made of short sentence-runs that all share
# identical run formatting (rFonts/color/sz) within each "block" (a block
# being the text between hard line breaks / lastRenderedPageBreak marks).
# Word (and this COM runtime) automatically coalesces adjacent runs that
# share identical formatting whenever any one of them is edited, so each
# block is replaced here as a single Find/Replace covering its whole text.

$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute could not find text: $old"
    }
}

# Title
$oldTitle = "The Enigmatic Cosmos: Unveiling Dark Energy's Secrets"
$newTitle = "The Profound Evolution of Artistic Expression: Embracing Diversity and Contemporary Forms"
ReplaceText $oldTitle $newTitle

# Author name ("Dr" + "." + " Elara Vanderstelt" collapses into one run)
ReplaceText "Dr. Elara Vanderstelt" "Isabella Kingsley"

# Email address - first run's text changes; the unchanged ". org" tail
# merges right onto it, producing "isabella.kingsley@educationalhaven.org"
ReplaceText "evanderstelt@cosmicmysteries" "isabella.kingsley@educationalhaven"

# Body paragraph, first block (before the first line break)
$old1 = "Across the vast expanse of the cosmos, amidst shimmering constellations and radiant nebulae, " + `
        "lies a profound enigma that captivates and perplexes scientists: dark energy." + `
        " This elusive force permeates the universe, driving its accelerated expansion at an ever-increasing rate." + `
        " Unraveling the secrets of dark energy holds the key to understanding the ultimate fate of our universe, yet its true nature remains shrouded in mystery." + `
        " Embark on a cosmic journey as we delve into the depths of this enigmatic phenomenon, exploring its profound implications for our comprehension of space, time, and the fundamental forces that govern our existence."
$new1 = "In the ever-evolving realm of artistic expression, the landscape of creativity has been dramatically transformed over time." + `
        " From the dawn of humanity, art has served as a profound and multifaceted means of communication, allowing individuals to connect with their innermost thoughts, experiences, and surroundings." + `
        " In this essay, we will explore the captivating journey of artistic evolution, delving into its diverse manifestations and embracing the contemporary forms that continue to push the boundaries of creativity."
ReplaceText $old1 $new1

# Body paragraph, second block (between the two line breaks)
$old2 = "The existence of dark energy was first hinted at in the late 1990s when astronomers observed that the expansion of the universe was not decelerating as expected, but rather accelerating." + `
        " This unexpected discovery challenged prevailing theories and ignited a scientific quest to identify the mysterious force responsible." + `
        " Enter dark energy, a hypothetical form of energy permeating the entire universe and causing its expansion to accelerate." + `
        " But what is this enigmatic entity? Is it a cosmological constant, a manifestation of vacuum energy, or something else entirely?"
$new2 = "Throughout history, art has undergone remarkable shifts, reflecting the cultural, social, and technological changes that have shaped human civilization." + `
        " The Renaissance witnessed a resurgence of classical ideals, while the Baroque period exuded grandeur and flamboyance." + `
        " Modernism shattered traditional norms, giving rise to bold abstraction and unconventional artistic expressions." + `
        " These eras of artistic evolution have left an indelible mark on our understanding of beauty, creativity, and the human condition."
ReplaceText $old2 $new2

# Body paragraph, third block (after the second line break)
$old3 = "As scientists grapple with these perplexing questions, the study of dark energy has illuminated other cosmic mysteries." + `
        " For instance, it has provided insights into the geometry of the universe, suggesting it may be flat or even slightly curved." + `
        " Moreover, dark energy's influence on the cosmic microwave background radiation has helped refine our understanding of the universe's composition and evolution." + `
        " While these discoveries have expanded our knowledge, they have also deepened the enigma surrounding dark energy, beckoning us to unravel its true identity and uncover its profound implications for the cosmos."
$new3 = "As we venture into the contemporary art scene, we encounter a kaleidoscope of diverse practices that challenge conventional notions of artistic expression." + `
        " Street art adorns urban landscapes, conveying powerful messages of social commentary and political resistance." + `
        " Digital art harnesses the boundless possibilities of technology, blurring the lines between reality and the virtual realm." + `
        " Performance art pushes the boundaries of artistic engagement, inviting audiences to become active participants in the creative process." + `
        " These contemporary forms reflect the dynamic and interconnected world we inhabit, showcasing the limitless potential of artistic innovation."
ReplaceText $old3 $new3

# Summary paragraph, first block (before the lastRenderedPageBreak run)
$old4 = "Dark energy, a mysterious and enigmatic force, permeates the universe, driving its accelerated expansion." + `
        " Its existence and nature remain shrouded in uncertainty, challenging prevailing "
$new4 = "The journey of artistic evolution has been a captivating tapestry of creativity, reflecting the cultural, social, and technological transformations that have shaped human civilization." + `
        " From the Renaissance to the Baroque period, from Modernism to the diverse contemporary forms, art "
ReplaceText $old4 $new4

# Summary paragraph, the lastRenderedPageBreak run itself (stays its own run)
ReplaceText "theories and captivating scientists worldwide" `
            "has served as a profound means of communication, allowing individuals to connect with their innermost thoughts, experiences, and surroundings"

# Summary paragraph, final block (after the lastRenderedPageBreak run)
$old5 = "." + `
        " The quest to understand dark energy has illuminated other cosmic mysteries, providing insights into the geometry of the universe and the composition of the cosmic microwave background radiation." + `
        " However, these discoveries have also intensified the enigma, urging us to unveil the true identity of dark energy and fathom its profound implications, not only for our understanding of the cosmos but also for the ultimate fate of our universe."
$new5 = "." + `
        " As we continue to embrace the ever-evolving landscape of artistic expression, we can appreciate the boundless potential of human imagination and creativity."
ReplaceText $old5 $new5

# A trailing empty paragraph is appended at the very end of the document.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
